$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value reads as a plain number (e.g. "212.61").
# Force text formatting first so Excel stores them as strings, matching the
# original inline-string (price-as-text) representation instead of
# converting them to numeric cells.
$numericLookingCells = @("D5","D8","D14","D16","D18","D20","D23","D25","D29","D33","D38","D42","D43","D44","D46","D48","D49","D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.550.03"
$ws.Range("E2").Value = "  -0.12%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.646.75"
$ws.Range("E3").Value = "  -1.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "212.61"
$ws.Range("E5").Value = "  -1.32%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +4.14%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.09%  "

# Row 8 - Solana
$ws.Range("D8").Value = "23.61"
$ws.Range("E8").Value = "  -2.63%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.86%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.32%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.54%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.879.47"
$ws.Range("E12").Value = "  -1.10%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.649.72"
$ws.Range("E13").Value = "  +1.20%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "0.593"
$ws.Range("E14").Value = "  +4.07%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  -2.15%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "64.55"
$ws.Range("E16").Value = "  -2.56%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.518.53"
$ws.Range("E17").Value = "  -0.20%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "232.05"
$ws.Range("E18").Value = "  -3.71%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  -0.90%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "7.56"
$ws.Range("E20").Value = "  -1.37%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.10%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -3.62%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "9.80"
$ws.Range("E23").Value = "  +4.48%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -1.45%  "

# Row 25 - Monero
$ws.Range("D25").Value = "148.66"
$ws.Range("E25").Value = "  +1.66%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -2.87%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +1.53%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  -0.12%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "15.64"
$ws.Range("E29").Value = "  -4.43%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -2.56%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -3.29%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.88%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "3.16"
$ws.Range("E33").Value = "  +1.07%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.424.50"
$ws.Range("E34").Value = "  -2.41%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +0.81%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +0.18%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -0.87%  "

# Row 38 - ARBITRUM
$ws.Range("D38").Value = "0.890"
$ws.Range("E38").Value = "  -4.32%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -3.27%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  -0.19%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  -0.03%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "0.819"
$ws.Range("E42").Value = "  +3.21%  "

# Row 43 - was FraxShare, becomes mCoin
$ws.Range("B43").Value = "mCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D43").Value = "2.47"
$ws.Range("E43").Value = "  -1.52%  "

# Row 44 - was mCoin, becomes FraxShare
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.54"
$ws.Range("E44").Value = "  +2.19%  "

# Row 45 - MXToken
$ws.Range("E45").Value = "  +0.94%  "

# Row 46 - Aave
$ws.Range("D46").Value = "65.16"
$ws.Range("E46").Value = "  -6.89%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.788.73"
$ws.Range("E47").Value = "  -1.07%  "

# Row 48 - RenderToken
$ws.Range("D48").Value = "1.69"
$ws.Range("E48").Value = "  -1.96%  "

# Row 49 - Quant
$ws.Range("D49").Value = "88.14"
$ws.Range("E49").Value = "  -0.77%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  -0.03%  "

# Row 51 - was EnergySwap, becomes Algorand
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0994"
$ws.Range("E51").Value = "  -3.41%  "
